$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state (post "Natmi following Dr Hou advice" re-run): the
# sending/target cluster grid grows from 3 target clusters (ECs, FAPs,
# sCs) to 4 (ECs, FAPs, M2, sCs), and gains a 4th sending cluster (M2)
# as well, giving a full 4x4 grid of 16 data rows (rows 2-17) for the
# Dhh -> Ptch2 ligand-receptor pair.
#
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity,
# T Edge total expression derived specificity.

$data = @(
  @("ECs","Dhh","Ptch2","ECs",3,1,1.569424333333333,4.708273,0.2941318537634731,0.2941318537634731,3,1,4.293676666666666,12.88103,0.2876481693908581,0.2876481693908581,6.738600640132222,60.64740576119,0.08460648929460261,0.08460648929460263),
  @("ECs","Dhh","Ptch2","FAPs",3,1,1.569424333333333,4.708273,0.2941318537634731,0.2941318537634731,3,1,8.749439000000001,26.248317,0.5861550151378376,0.5861550151378377,13.73158246961567,123.584242226541,0.1724068611952488,0.1724068611952489),
  @("ECs","Dhh","Ptch2","M2",3,1,1.569424333333333,4.708273,0.2941318537634731,0.2941318537634731,3,1,0.6753636666666667,2.026091,0.04524493516196244,0.04524493516196246,1.059932172315889,9.539389550843001,0.01330797665259616,0.01330797665259616),
  @("ECs","Dhh","Ptch2","sCs",3,1,1.569424333333333,4.708273,0.2941318537634731,0.2941318537634731,3,1,1.208355333333333,3.625066,0.08095188030934176,0.08095188030934178,1.896422263446444,17.067800371018,0.02381052662102549,0.02381052662102549),
  @("FAPs","Dhh","Ptch2","ECs",3,1,2.444496,7.333488,0.4581324022613356,0.4581324022613355,3,1,4.293676666666666,12.88103,0.2876481693908581,0.2876481693908581,10.49587543696,94.46287893264,0.1317809468491094,0.1317809468491094),
  @("FAPs","Dhh","Ptch2","FAPs",3,1,2.444496,7.333488,0.4581324022613356,0.4581324022613355,3,1,8.749439000000001,26.248317,0.5861550151378376,0.5861550151378377,21.387968637744,192.491717739696,0.2685366051826271,0.2685366051826271),
  @("FAPs","Dhh","Ptch2","M2",3,1,2.444496,7.333488,0.4581324022613356,0.4581324022613355,3,1,0.6753636666666667,2.026091,0.04524493516196244,0.04524493516196246,1.650923781712,14.858314035408,0.02072817083590823,0.02072817083590823),
  @("FAPs","Dhh","Ptch2","sCs",3,1,2.444496,7.333488,0.4581324022613356,0.4581324022613355,3,1,1.208355333333333,3.625066,0.08095188030934176,0.08095188030934178,2.953819778912,26.584378010208,0.03708667939369085,0.03708667939369086),
  @("M2","Dhh","Ptch2","ECs",2,0.6666666666666666,0.165334,0.4960020000000001,0.03098588117774611,0.03098588117774611,3,1,4.293676666666666,12.88103,0.2876481693908581,0.2876481693908581,0.7098907380066667,6.389016642060001,0.008913031997741313,0.008913031997741315),
  @("M2","Dhh","Ptch2","FAPs",2,0.6666666666666666,0.165334,0.4960020000000001,0.03098588117774611,0.03098588117774611,3,1,8.749439000000001,26.248317,0.5861550151378376,0.5861550151378377,1.446579747626,13.019217728634,0.01816252965080101,0.01816252965080101),
  @("M2","Dhh","Ptch2","M2",2,0.6666666666666666,0.165334,0.4960020000000001,0.03098588117774611,0.03098588117774611,3,1,0.6753636666666667,2.026091,0.04524493516196244,0.04524493516196246,0.1116605764646667,1.004945188182,0.001401954184823395,0.001401954184823395),
  @("M2","Dhh","Ptch2","sCs",2,0.6666666666666666,0.165334,0.4960020000000001,0.03098588117774611,0.03098588117774611,3,1,1.208355333333333,3.625066,0.08095188030934176,0.08095188030934178,0.1997822206813334,1.798039986132,0.002508365344380389,0.002508365344380389),
  @("sCs","Dhh","Ptch2","ECs",3,1,1.156530666666667,3.469592,0.2167498627974453,0.2167498627974453,3,1,4.293676666666666,12.88103,0.2876481693908581,0.2876481693908581,4.965768737751111,44.69191863976,0.06234770124940479,0.0623477012494048),
  @("sCs","Dhh","Ptch2","FAPs",3,1,1.156530666666667,3.469592,0.2167498627974453,0.2167498627974453,3,1,8.749439000000001,26.248317,0.5861550151378376,0.5861550151378377,10.11899451962933,91.070950676664,0.1270490191091608,0.1270490191091608),
  @("sCs","Dhh","Ptch2","M2",3,1,1.156530666666667,3.469592,0.2167498627974453,0.2167498627974453,3,1,0.6753636666666667,2.026091,0.04524493516196244,0.04524493516196246,0.7810787916524444,7.029709124872,0.009806833488634668,0.00980683348863467),
  @("sCs","Dhh","Ptch2","sCs",3,1,1.156530666666667,3.469592,0.2167498627974453,0.2167498627974453,3,1,1.208355333333333,3.625066,0.08095188030934176,0.08095188030934178,1.397499999230222,12.577499993072,0.01754630895024504,0.01754630895024504)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
